$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961541110818"
$ws1.Range("B2").Value = "go_stims-16509961540790472.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961540950804.csv"
$ws1.Range("B4").Value = "go_stims-16509961540950804.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961541110818.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961557430844"
$ws2.Range("B2").Value = "OB-16509961546950452.csv"
$ws2.Range("B3").Value = "TB-165099615571908.csv"
$ws2.Range("B4").Value = "OB-1650996155463039.csv"
$ws2.Range("B5").Value = "ZB-match_4-16509961545430434.csv"
$ws2.Range("B6").Value = "TB-1650996155687077.csv"
$ws2.Range("B7").Value = "TB-1650996155479126.csv"
$ws2.Range("B8").Value = "OB-16509961551990733.csv"
$ws2.Range("B9").Value = "ZB-match_6-16509961542230816.csv"
$ws2.Range("B10").Value = "ZB-match_7-16509961544950445.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961557430844"
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961557910767"
$ws4.Range("B2").Value = "MM_stims-16509961557590482.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961557430844.csv"
$ws4.Range("B4").Value = "MM_stims-16509961557750418.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961557590482.csv"
$ws4.Range("B6").Value = "MM_stims-16509961557910767.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961557750418.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509961558552272"
$ws5.Range("B2").Value = "SAT_stims-16509961557910767.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961558071926.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961558392272.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650996155823231.csv"
